# Update "Generate Report for Handback" timestamps on the zh-cn and de-de sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 13:04:58"
$wsZhCn.Range("H2").Value = "2016-03-13 13:05:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 13:05:02"
$wsDeDe.Range("H2").Value = "2016-03-13 13:05:21"
